$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stats for row 6 (new match results added in)
$ws.Range("B6").Value = 29
$ws.Range("C6").Value = 13
$ws.Range("D6").Value = 25
$ws.Range("E6").Value = 132
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = 3438
$ws.Range("H6").Value = 44

# Update stats for row 7
$ws.Range("B7").Value = 23
$ws.Range("C7").Value = 11
$ws.Range("D7").Value = 30
$ws.Range("F7").Value = 4
$ws.Range("G7").Value = 2990
$ws.Range("H7").Value = 44

# Update stats for row 12
$ws.Range("B12").Value = 30
$ws.Range("C12").Value = 22
$ws.Range("D12").Value = 28
$ws.Range("F12").Value = 13
$ws.Range("G12").Value = 3678
$ws.Range("H12").Value = 39

# Update stats for row 13
$ws.Range("B13").Value = 7
$ws.Range("D13").Value = 29
$ws.Range("F13").Value = 3
$ws.Range("G13").Value = 1114
$ws.Range("H13").Value = 39

# Move the active selection to reflect where the author left off editing
$ws.Range("R14").Select()
